$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.209.22"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "2.087.44"
$ws.Range("E3").Value = "  +3.26%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'250.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.96%  "
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("D8").Value = "'54.56"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +24.93%  "
$ws.Range("D9").Value = "'61.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("D10").Value = "'0.377"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.10%  "
$ws.Range("D11").Value = "'0.0746"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.33%  "
$ws.Range("E12").Value = "  +8.28%  "
$ws.Range("D13").Value = "'15.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.11%  "
$ws.Range("D14").Value = "2.392.00"
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("D15").Value = "'0.832"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.06%  "
$ws.Range("D16").Value = "2.085.12"
$ws.Range("E16").Value = "  +3.08%  "
$ws.Range("D17").Value = "'5.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.45%  "
$ws.Range("D18").Value = "37.138.38"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").Value = "'72.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("D20").Value = "'14.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +17.09%  "
$ws.Range("D21").Value = "0.0₃0848"
$ws.Range("E21").Value = "  +5.65%  "
$ws.Range("D22").Value = "'240.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("D23").Value = "'5.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.88%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'2.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").Value = "'172.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.54%  "
$ws.Range("D27").Value = "'9.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.86%  "
$ws.Range("E28").Value = "  +5.35%  "
$ws.Range("D29").Value = "'2.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.66%  "
$ws.Range("D30").Value = "'0.124"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.07%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +32.10%  "
$ws.Range("B32").Value = "Gas"
$ws.Range("C32").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D32").Value = "'22.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.66%  "
$ws.Range("D33").Value = "'4.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.22%  "
$ws.Range("D34").Value = "'0.0620"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.73%  "
$ws.Range("D35").Value = "'0.0902"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.46%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "'4.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.07%  "
$ws.Range("E38").Value = "  -3.34%  "
$ws.Range("D39").Value = "'2.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.10%  "
$ws.Range("E40").Value = "  +3.75%  "
$ws.Range("D41").Value = "'4.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +152.17%  "
$ws.Range("E42").Value = "  +19.04%  "
$ws.Range("E43").Value = "  +7.53%  "
$ws.Range("D44").Value = "'1.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.48%  "
$ws.Range("D45").Value = "'98.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.14%  "
$ws.Range("D46").Value = "'0.0949"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +16.84%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").Value = "1.325.29"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("E49").Value = "  +5.34%  "
$ws.Range("E50").Value = "  +9.26%  "
$ws.Range("D51").Value = "'6.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +15.73%  "
